$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.788.48"
$ws.Range("E2").Value = "  -8.86%  "

# Row 3
$ws.Range("D3").Value = "2.892.08"
$ws.Range("E3").Value = "  -8.66%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.00"
$ws.Range("E5").Value = "  -8.99%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "120.83"
$ws.Range("E6").Value = "  -11.06%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").Value = "2.888.18"
$ws.Range("E8").Value = "  -8.82%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -4.09%  "

# Row 10
$ws.Range("E10").Value = "  -12.03%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.81"
$ws.Range("E11").Value = "  -10.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.430"
$ws.Range("E12").Value = "  -5.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000212"
$ws.Range("E13").Value = "  -11.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.27"
$ws.Range("E14").Value = "  -10.34%  "

# Row 15
$ws.Range("E15").Value = "  -1.52%  "

# Row 16
$ws.Range("D16").Value = "3.358.92"
$ws.Range("E16").Value = "  -9.01%  "

# Row 17
$ws.Range("D17").Value = "2.892.30"
$ws.Range("E17").Value = "  -8.66%  "

# Row 18
$ws.Range("D18").Value = "57.897.18"
$ws.Range("E18").Value = "  -8.71%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.32"
$ws.Range("E19").Value = "  -3.74%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "414.99"
$ws.Range("E20").Value = "  -10.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.74"
$ws.Range("E21").Value = "  -8.76%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.648"
$ws.Range("E22").Value = "  -6.98%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.77"
$ws.Range("E23").Value = "  -11.47%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.51"
$ws.Range("E24").Value = "  -6.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "76.61"
$ws.Range("E25").Value = "  -7.87%  "

# Row 26
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("E27").Value = "  -0.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.90"
$ws.Range("E29").Value = "  -8.85%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("E30").Value = "  -9.22%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "24.54"
$ws.Range("E31").Value = "  -9.37%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.92"
$ws.Range("E32").Value = "  -12.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0931"
$ws.Range("E33").Value = "  -7.25%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.34"
$ws.Range("E34").Value = "  -9.41%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "48.53"
$ws.Range("E35").Value = "  -5.19%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.885"
$ws.Range("E36").Value = "  -13.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.97"
$ws.Range("E37").Value = "  -18.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.32"
$ws.Range("E38").Value = "  +2.38%  "

# Row 39
$ws.Range("D39").Value = "0.0₃0621"
$ws.Range("E39").Value = "  -14.82%  "

# Row 40
$ws.Range("E40").Value = "  -12.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.104"
$ws.Range("E41").Value = "  -7.82%  "

# Row 42
$ws.Range("D42").Value = "2.607.22"
$ws.Range("E42").Value = "  -6.50%  "

# Row 43
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "349.40"
$ws.Range("E43").Value = "  -10.63%  "

# Row 44
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("E45").Value = "  -10.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "118.31"
$ws.Range("E46").Value = "  -6.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.226"
$ws.Range("E47").Value = "  -9.41%  "

# Row 48
$ws.Range("E48").Value = "  -5.51%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.92"
$ws.Range("E49").Value = "  -9.06%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.50"
$ws.Range("E50").Value = "  -10.18%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.93"
$ws.Range("E51").Value = "  -10.65%  "
